$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (AD, AE, AF): Wins / Losses / Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header formatting (bold, centered, bordered) from an existing header cell
# onto the new header cells so they match the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the season record (Wins / Losses / Ties) for every data row (2 through 50)
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 63
    $ws.Cells.Item($row, 31).Value = 99
    $ws.Cells.Item($row, 32).Value = 0
}
